$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph, before "Gameplay Features".
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$null = $titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaXml = '<w:p ' + $wNs + '>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
    '<w:r><w:t>: Read our review of Dead or Alive slot game by NetEnt and play for free. Enjoy the immersive Western-themed atmosphere and Sticky Win feature for big wins.</w:t></w:r>' +
    '</w:p>'
$null = $metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Dead or Alive Slot Game for
#    Free" paragraph near the end of the document, and replace the
#    text of the following italic paragraph with the new image
#    description prompt.
# ------------------------------------------------------------------
$dupTitleText = "Play Dead or Alive Slot Game for Free"
foreach ($p in $d.Paragraphs) {
    $sameText = $p.Range.Text.TrimEnd([char]13, [char]7) -eq $dupTitleText
    if ($sameText -and $p.Style.NameLocal -ne "Heading 1") {
        $null = $p.Range.Delete()
        break
    }
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newImagePrompt = "Create a cartoon-style feature image that captures the essence of Dead or Alive game by featuring a happy Maya warrior with glasses. The image can have the warrior holding a smoking gun or standing in front of a Western saloon. The colors used in the image should be dark, representing the mood of the game while also highlighting the warrior's colorful attire. The warrior should have a big smile on their face, creating a contrast with the game's theme, showing that the game is enjoyable and fun to play. The image should also contain the game's name, `"Dead or Alive`" prominently displayed in bold and colorful fonts."

$imageXml = '<w:p ' + $wNs + '>' +
    '<w:r/>' +
    '<w:r><w:rPr><w:i/></w:rPr><w:t>' + $newImagePrompt + '</w:t></w:r>' +
    '</w:p>'
$null = $lastPara.Range.InsertXML($imageXml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
